$d = $word.ActiveDocument

# --- Locate the target paragraph ("Baz chan" + bookmark + "ges") ---
# Scope every Find to this paragraph's range so we never touch the word
# "changes" that already appears earlier in the document's second paragraph.
$target = $d.Paragraphs(5).Range

# 1) Replace "Baz chan" with the new author line.
[void]$target.Find.Execute("Baz chan", $true, $false, $false, $false, $false, `
                      $true, 1, $false, `
                      "A23lawson " + [char]0x2013 + " Ashley Lawson", 2)

# 2) Temporarily remove the _GoBack bookmark - it sits right between the two
#    runs and would otherwise "stick" to whichever side we don't want once we
#    split the paragraph.
$bm = $d.Bookmarks("_GoBack")
$splitPos = $bm.Range.End
$bm.Delete()

# 3) Split the paragraph in two at the old bookmark location by inserting a
#    paragraph mark there.
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertAfter([char]13)

# 4) The tail half is now its own paragraph containing just "ges" - turn it
#    into the full second line of text.
$tail = $d.Paragraphs(6).Range
[void]$tail.Find.Execute("ges", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "Baz changes", 2)

# 5) Re-create the _GoBack bookmark at the end of the first (new) paragraph,
#    after its text, so it again reads "<text><bookmark/></p>".
#    Inserting a bookmark with a zero-length range sitting exactly on a
#    paragraph's end boundary is unreliable, so append a one-character
#    placeholder first, anchor the bookmark there, then delete the
#    placeholder - leaving the bookmark correctly positioned.
$head = $d.Paragraphs(5).Range
$bmPos = $head.End - 1
$placeholder = $d.Range($bmPos, $bmPos)
$placeholder.InsertAfter("X")

$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))

$bmNow = $d.Bookmarks("_GoBack")
$afterBm = $bmNow.Range.End
$d.Range($afterBm, $afterBm + 1).Delete()
